# Updated cryptos list on Sat Apr 22 21:40:54 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to stay plain text (source values like
# "1.003" or "27.736.02" must not be auto-coerced into numbers/dates).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.736.02"
$ws.Range("E2").Value = "  +1.21%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.876.48"
$ws.Range("E3").Value = "  +1.39%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5 - BNB
$ws.Range("D5").Value = "333.19"
$ws.Range("E5").Value = "  +3.65%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.16%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4742"
$ws.Range("E7").Value = "  +6.35%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3964"
$ws.Range("E8").Value = "  +3.43%  "

# Row 9 - OKB
$ws.Range("D9").Value = "47.80"
$ws.Range("E9").Value = "  -2.87%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.08042"
$ws.Range("E10").Value = "  +2.74%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "1.025"
$ws.Range("E11").Value = "  +0.98%  "

# Row 12 - Solana
$ws.Range("D12").Value = "21.89"
$ws.Range("E12").Value = "  +2.22%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.876.91"
$ws.Range("E13").Value = "  +1.45%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "5.964"
$ws.Range("E14").Value = "  +2.15%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "7.165"
$ws.Range("E15").Value = "  +0.97%  "

# Row 16 - BinanceUSD (D unchanged)
$ws.Range("E16").Value = "  +0.46%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.00001049"
$ws.Range("E17").Value = "  +2.27%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "87.26"
$ws.Range("E18").Value = "  +2.39%  "

# Row 19 - TRON (E unchanged)
$ws.Range("D19").Value = "0.06636"

# Row 20 - Avalanche
$ws.Range("D20").Value = "17.32"
$ws.Range("E20").Value = "  +2.31%  "

# Row 21 - Dai (D unchanged)
$ws.Range("E21").Value = "  +0.14%  "

# Row 22 - WrappedBTC
$ws.Range("D22").Value = "27.724.16"
$ws.Range("E22").Value = "  +1.18%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.497"
$ws.Range("E23").Value = "  +0.57%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  +2.59%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.300"
$ws.Range("E25").Value = "  +1.80%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").Value = "2.099.89"
$ws.Range("E26").Value = "  +1.58%  "

# Row 27 - Monero
$ws.Range("D27").Value = "156.69"
$ws.Range("E27").Value = "  +3.55%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "20.22"
$ws.Range("E28").Value = "  +4.71%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "2.104"
$ws.Range("E29").Value = "  +2.95%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "5.585"
$ws.Range("E30").Value = "  +2.16%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "122.55"
$ws.Range("E31").Value = "  +2.11%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "0.9708"
$ws.Range("E32").Value = "  +4.85%  "

# Row 33 - Stellar (D unchanged)
$ws.Range("E33").Value = "  +2.74%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "1.457"
$ws.Range("E34").Value = "  -1.04%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "3.637"
$ws.Range("E35").Value = "  +0.33%  "

# Row 36 - Filecoin
$ws.Range("D36").Value = "5.305"
$ws.Range("E36").Value = "  +1.79%  "

# Row 37 - Hedera
$ws.Range("D37").Value = "0.06115"
$ws.Range("E37").Value = "  +2.96%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.02265"
$ws.Range("E38").Value = "  +2.22%  "

# Row 39 - TrustWalletToken
$ws.Range("D39").Value = "1.229"
$ws.Range("E39").Value = "  +2.21%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "8.201"
$ws.Range("E40").Value = "  -1.36%  "

# Row 41 - Frax (D unchanged)
$ws.Range("E41").Value = "  +0.22%  "

# Row 42 - TheSandbox
$ws.Range("D42").Value = "0.5993"
$ws.Range("E42").Value = "  +1.49%  "

# Row 43 - Algorand
$ws.Range("D43").Value = "0.1913"
$ws.Range("E43").Value = "  +3.63%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "10.29"
$ws.Range("E44").Value = "  +0.45%  "

# Row 45 & 46 swap: WEMIXTOKEN and Decentraland traded places in the ranking
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5723"
$ws.Range("E45").Value = "  +1.09%  "

$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.252"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47 - EnergySwap (D unchanged)
$ws.Range("E47").Value = "  +1.56%  "

# Row 48 - PancakeSwap
$ws.Range("D48").Value = "3.411"
$ws.Range("E48").Value = "  +1.58%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "1.938"
$ws.Range("E49").Value = "  +1.16%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.06811"
$ws.Range("E50").Value = "  -0.64%  "

# Row 51 - Quant
$ws.Range("D51").Value = "112.86"
$ws.Range("E51").Value = "  +4.51%  "

# Restore the original (General) number format now that the text values
# are locked in, so cell formatting is unchanged from the source file.
$ws.Range("D2:E51").NumberFormat = "General"
